$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.247.97"
$ws.Cells.Item(2, 5).Value = "  -0.11%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.591.84"
$ws.Cells.Item(3, 5).Value = "  +0.20%  "

$ws.Cells.Item(4, 5).Value = "  -0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.43"
$ws.Cells.Item(5, 5).Value = "  +1.03%  "

$ws.Cells.Item(6, 5).Value = "  -0.57%  "

$ws.Cells.Item(7, 5).Value = "  -0.11%  "

$ws.Cells.Item(8, 5).Value = "  -0.22%  "

$ws.Cells.Item(9, 5).Value = "  -0.46%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.07"
$ws.Cells.Item(10, 5).Value = "  -1.46%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0850"
$ws.Cells.Item(11, 5).Value = "  +0.57%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.814.79"
$ws.Cells.Item(12, 5).Value = "  +0.13%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.593.59"
$ws.Cells.Item(13, 5).Value = "  -0.11%  "

$ws.Cells.Item(14, 5).Value = "  -1.83%  "

$ws.Cells.Item(15, 5).Value = "  -2.09%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "63.87"
$ws.Cells.Item(16, 5).Value = "  -0.89%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "26.226.58"
$ws.Cells.Item(17, 5).Value = "  -0.25%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0₃0727"
$ws.Cells.Item(18, 5).Value = "  -0.62%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "216.34"
$ws.Cells.Item(19, 5).Value = "  +2.58%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.30"
$ws.Cells.Item(20, 5).Value = "  -2.43%  "

$ws.Cells.Item(21, 5).Value = "  -0.01%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.29"
$ws.Cells.Item(22, 5).Value = "  +0.14%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.07"
$ws.Cells.Item(23, 5).Value = "  +0.98%  "

$ws.Cells.Item(24, 5).Value = "  -0.72%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.56"
$ws.Cells.Item(25, 5).Value = "  +0.35%  "

$ws.Cells.Item(26, 5).Value = "  -0.13%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.98"
$ws.Cells.Item(27, 5).Value = "  -1.16%  "

$ws.Cells.Item(28, 5).Value = "  -0.69%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.14"
$ws.Cells.Item(29, 5).Value = "  -0.50%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0492"
$ws.Cells.Item(30, 5).Value = "  -2.16%  "

$ws.Cells.Item(31, 5).Value = "  +0.51%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.19"
$ws.Cells.Item(32, 5).Value = "  -0.53%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.428.95"
$ws.Cells.Item(33, 5).Value = "  +8.13%  "

$ws.Cells.Item(34, 5).Value = "  -0.74%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.43"
$ws.Cells.Item(35, 5).Value = "  -0.61%  "

$ws.Cells.Item(36, 5).Value = "  -0.41%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.584"
$ws.Cells.Item(37, 5).Value = "  -3.07%  "

$ws.Cells.Item(38, 5).Value = "  -0.96%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.827"
$ws.Cells.Item(39, 5).Value = "  +1.94%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "5.92"
$ws.Cells.Item(40, 5).Value = "  +4.08%  "

$ws.Cells.Item(41, 5).Value = "  -0.10%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.988"
$ws.Cells.Item(42, 5).Value = "  -2.80%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.768"
$ws.Cells.Item(43, 5).Value = "  +0.36%  "

$ws.Cells.Item(44, 5).Value = "  +0.11%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.726.98"
$ws.Cells.Item(45, 5).Value = "  +0.10%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "61.12"
$ws.Cells.Item(46, 5).Value = "  -1.26%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "86.96"
$ws.Cells.Item(47, 5).Value = "  -1.02%  "

$ws.Cells.Item(48, 5).Value = "  +0.57%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0502"
$ws.Cells.Item(49, 5).Value = "  -0.54%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0954"
$ws.Cells.Item(50, 5).Value = "  -2.08%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.999"
$ws.Cells.Item(51, 5).Value = "  -0.13%  "
